# Append " (Changed main)" to the end of the first paragraph's text,
# as three distinct runs: " (", "Changed main", ")".
#
# Word normally coalesces successive InsertAfter calls with identical
# run formatting into a single run. Briefly enabling TrackRevisions
# around the inserts keeps each insertion as its own run (tracked
# insertions never merge with neighbouring text), then AcceptAllRevisions
# folds the tracked markup back into plain runs without re-merging them.
$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$r = $p1.Range
# Paragraph.Range includes the trailing paragraph mark; back off one
# character so the collapsed insertion point stays inside paragraph 1.
$r.End = $r.End - 1

$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true

$r.Collapse(0)
$r.InsertAfter(" (")

$r.Collapse(0)
$r.InsertAfter("Changed main")

$r.Collapse(0)
$r.InsertAfter(")")

$d.TrackRevisions = $wasTracking
$d.AcceptAllRevisions()
